$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the header in D1 from "2ndEmail" to "Contact2"
$ws.Range("D1").Value = "Contact2"

# Add a new row of data (row 20): First Name = arnold, Email = akillingbeck@gmail.com
$ws.Range("C20").Value = "akillingbeck@gmail.com"
$ws.Range("A20").Value = "arnold"

# Update selection to match the saved state
$ws.Range("A20").Select()
